$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# Text-like columns (Date, Time, Weekday, Week) must stay as text, not be
# coerced into dates/numbers by Excel's automatic type inference.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-02-27"

$ws.Range("B$row").NumberFormat = "@"
$ws.Range("B$row").Value = "08:56:47"

$ws.Range("C$row").Value = "Thursday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "08"

# Numeric columns
$ws.Range("E$row").Value = 131153
$ws.Range("F$row").Value = 142078
$ws.Range("G$row").Value = 172871
$ws.Range("H$row").Value = 160217
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 146543
$ws.Range("K$row").Value = -1
$ws.Range("L$row").Value = -1
$ws.Range("M$row").Value = 194084
$ws.Range("N$row").Value = 115483
$ws.Range("O$row").Value = 46810
$ws.Range("P$row").Value = 29557
$ws.Range("Q$row").Value = 69456
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 49444
$ws.Range("T$row").Value = -1
